$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -8.207999999999998
$ws.Range("B7").Value = 5.814
$ws.Range("A9").Value = -21.396
$ws.Range("B12").Value = 5.927
$ws.Range("A13").Value = -22.116
$ws.Range("B14").Value = 6.040000000000001
$ws.Range("D15").Value = -8.202000000000002
$ws.Range("A16").Value = -21.886
$ws.Range("A18").Value = -22.055
$ws.Range("B19").Value = 8.395999999999999
$ws.Range("A20").Value = -20.32
$ws.Range("A26").Value = -21.396
$ws.Range("B26").Value = 5.647
$ws.Range("A27").Value = -21.188
$ws.Range("B27").Value = 5.695000000000001
$ws.Range("D28").Value = -8.205
$ws.Range("A29").Value = -21.722
$ws.Range("B29").Value = 6.068
$ws.Range("D33").Value = -7.627
$ws.Range("A35").Value = -20.111
$ws.Range("D35").Value = -7.532000000000001
$ws.Range("A36").Value = -21.095
$ws.Range("B37").Value = 8.847000000000001
$ws.Range("B38").Value = 5.920000000000001
$ws.Range("D38").Value = -8.821
$ws.Range("D43").Value = -7.801
$ws.Range("D44").Value = -7.465999999999999
$ws.Range("A45").Value = -21.704
$ws.Range("D45").Value = -7.565
$ws.Range("B47").Value = 6.043000000000001
$ws.Range("D47").Value = -8.051
$ws.Range("B51").Value = 5.465
$ws.Range("D51").Value = -8.190000000000001
$ws.Range("B52").Value = 5.395
$ws.Range("D54").Value = -8.306999999999999
$ws.Range("A55").Value = -21.682
$ws.Range("B55").Value = 6.206
$ws.Range("A57").Value = -22.258
$ws.Range("D57").Value = -8.231
$ws.Range("D62").Value = -7.904000000000001
$ws.Range("D63").Value = -7.336999999999999
$ws.Range("D67").Value = -6.853999999999999
$ws.Range("A69").Value = -21.641
$ws.Range("B69").Value = 5.857
$ws.Range("B70").Value = 5.667
$ws.Range("D70").Value = -6.976000000000001
$ws.Range("A76").Value = -22.022
$ws.Range("B76").Value = 5.055
$ws.Range("A78").Value = -20.192
$ws.Range("B81").Value = 6.197000000000001
$ws.Range("D81").Value = -7.413999999999999
$ws.Range("A82").Value = -22.143
$ws.Range("A83").Value = -20.292
$ws.Range("B83").Value = 7.202
$ws.Range("D88").Value = -8.198
$ws.Range("A93").Value = -21.941
$ws.Range("B94").Value = 7.105
$ws.Range("D96").Value = -7.467000000000001
$ws.Range("A97").Value = -22.053
$ws.Range("D99").Value = -8.101000000000001
$ws.Range("B100").Value = 5.116999999999999
$ws.Range("B102").Value = 7.575
